$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (trial numbers) for columns B:E
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 (CON): B2 removed entirely, C2:E2 updated with new mean EMG values
$ws.Range("B2").ClearContents()
$ws.Range("C2").Value = 1.1770765782808947
$ws.Range("D2").Value = 0.30258517878326446
$ws.Range("E2").Value = 3.4592388228240232

# Row 3 (STR): B3:E3 updated with new mean EMG values
$ws.Range("B3").Value = 0.63775836925333151
$ws.Range("C3").Value = 1.8997263969175724
$ws.Range("D3").Value = 2.0522215683671039
$ws.Range("E3").Value = 3.2769210066465044

# Update selection to match the reduced highlighted range
$ws.Range("B1:E3").Select()
